# Automatische test-sync: 2025-07-23 22:56:50
#
# Adds Testmail #19 ("Hoe gaan jullie om met mijn persoonsgegevens?") to the
# Logs sheet as row 29, extends the related conditional formatting ranges,
# refreshes the Dashboard category/count table (re-sorted order + the new
# "Juridisch / Contract" category in row 12) and widens the chart's series
# references to include the new Dashboard row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append the new test-mail row (row 29)
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A29").Value = "Hoe gaan jullie om met mijn persoonsgegevens?"
$logs.Range("B29").Value = "mailmind.test@zohomail.eu"
$logs.Range("C29").Value = "Testmail #19: Hoe gaan jullie om met mijn persoonsgegevens?"
$logs.Range("D29").Value = "Juridisch / Contract"

$antwoord29 = "Beste afzender,`nDank u voor uw interesse in onze aanpak met betrekking tot persoonsgegevens. Bij ons bedrijf hechten we veel waarde aan de bescherming van persoonlijke informatie. Wij voldoen aan alle geldende wet- en regelgeving met betrekking tot gegevensbescherming, inclusief de Algemene Verordening Gegevensbescherming (AVG).`nMocht u nog specifieke vragen hebben over hoe wij omgaan met persoonsgegevens of als u meer details wenst te ontvangen, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam assistent]`nNederlandse e-mailassistent`nBedrijfsnaam"
$logs.Range("E29").Value = $antwoord29

$logs.Range("F29").Value = "2025-07-23 22:56:35"
$logs.Range("G29").Value = "Ja"
$logs.Range("H29").Value = "Nee"
$logs.Range("I29").Value = "Ja"
$logs.Range("J29").Value = "Nee"

# Drop the auto row-height bump the multi-line Antwoord text triggers, so
# row 29 stays on the sheet's default row height like every other row.
$logs.Rows.Item(29).AutoFit()

# Extend the conditional-formatting sqref ranges (D/G/H/I/J) from row 28 to
# the new row 29, one ModifyAppliesToRange per block (each block holds every
# cfRule sharing that sqref).
$logs.Range("D2:D28").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D29"))
$logs.Range("G2:G28").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G29"))
$logs.Range("H2:H28").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H29"))
$logs.Range("I2:I28").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I29"))
$logs.Range("J2:J28").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J29"))

# ---------------------------------------------------------------------
# 2) Dashboard sheet: refresh the Categorie/Aantal summary table
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

# Rows 7-9 are re-sorted (IT / Technisch probleem moves ahead of Factuur /
# Administratie and Offerte / Prijsaanvraag); counts stay at 1 each.
$dash.Range("A7").Value = "IT / Technisch probleem"
$dash.Range("A8").Value = "Factuur / Administratie"
$dash.Range("A9").Value = "Offerte / Prijsaanvraag"

# New row 12 for the category introduced by the new test mail.
$dash.Range("A12").Value = "Juridisch / Contract"
$dash.Range("B12").Value = 1

# ---------------------------------------------------------------------
# 3) Chart: widen the category/value series references to row 12
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$12"
$series.Values = "='Dashboard'!`$B`$2:`$B`$12"
